# Phân công công việc - restore the "main" branch content (revert of a revert):
#  - Reassign "Người phụ trách" (column C) to the original roster
#  - Add a "Thời gian" (column D) schedule for each task
#  - Re-center several columns that previously lacked horizontal centering

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---- Header row (A3:D3): add horizontal centering ----
$ws.Range("A3:D3").HorizontalAlignment = $xlCenter

# ---- Column B (Công việc): add horizontal centering, keep existing wrap/fill ----
$ws.Range("B4:B11").HorizontalAlignment = $xlCenter

# ---- Column C (Người phụ trách): new assignments ----
$ws.Range("C4").Value = "Công"
$ws.Range("C5").Value = "N Tiến"
$ws.Range("C6").Value = "Công"
$ws.Range("C7").Value = "X Tiến"
$ws.Range("C8").Value = "Duy"
$ws.Range("C9").Value = "Công"
$ws.Range("C10").Value = "Vân"

# C11 previously had no border/wrap (different style family than C4/C7) - pull the
# bordered+wrapped format from C7 (its row is the closest existing analogue), then
# set its value.
$fmtSrc = $ws.Range("C7")
$fmtSrc.Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = "Vân"
$excel.CutCopyMode = $false

# ---- Column D (Thời gian): new schedule, short-date number format + centering ----
$ws.Range("D4:D11").NumberFormat = "d-mmm"
$ws.Range("D4:D11").HorizontalAlignment = $xlCenter

$ws.Range("D4").Value = "21-22/10/2020"
$ws.Range("D5").Value = "23-24/10/2020"
$ws.Range("D6").Value = "23-26/10/2020"
$ws.Range("D7").Value = "23-26/10/2020"
$ws.Range("D8").Value = "27-1/11/2020"
$ws.Range("D9").Value = "03-07/11/2020"
$ws.Range("D10").Value = "08-015/11/2020"
$ws.Range("D11").Value = "16-17/11/2020"

# D11 (and D7/D8 already wrap from the row's default format) must NOT wrap in the
# target layout - only D7/D8 keep wrap.
$ws.Range("D11").WrapText = $false

# ---- Selection moves to E4 ----
$ws.Range("E4").Select()
